$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# D2: "Nadeen Amr" -> "Nadeen" (frees up the now-unused "Nadeen Amr" shared string)
$ws.Range("D2").Value = "Nadeen"

# Row 18 / row 19 reviewer swap
$ws.Range("F18").Value = "Abougabal"
$ws.Range("E19").Value = "Naka"

# Move the active selection to F17 (matches the saved sheet view)
$ws.Range("F17").Select() | Out-Null
